$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 6699.4
$ws.Range("I43").Value = 5050
$ws.Range("J43").Value = 7799
$ws.Range("K43").Value = 5050
$ws.Range("L43").Value = 7799
$ws.Range("M43").Value = -4981
$ws.Range("N43").Value = -7937

# Row 94
$ws.Range("H94").Value = 601.7778
$ws.Range("I94").Value = 601.7778
$ws.Range("K94").Value = 601.7778
$ws.Range("M94").Value = -150.7778

# Row 137
$ws.Range("H137").Value = 2551.9473
$ws.Range("I137").Value = 1953.909
$ws.Range("J137").Value = 3374.25
$ws.Range("K137").Value = 5861.727000000001
$ws.Range("L137").Value = 10122.75
$ws.Range("M137").Value = -3311.727000000001
$ws.Range("N137").Value = -15222.75

# Row 138
$ws.Range("H138").Value = 2643.9333
$ws.Range("I138").Value = 1430.4138
$ws.Range("J138").Value = 3779.1614
$ws.Range("K138").Value = 4291.2414
$ws.Range("L138").Value = 11337.4842
$ws.Range("M138").Value = 848.7586000000001
$ws.Range("N138").Value = -21617.4842

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4043.8667
$ws.Range("I32").Value = 3624.4092
$ws.Range("K32").Value = 3624.4092
$ws.Range("M32").Value = -3337.4092

# Row 61
$ws.Range("H61").Value = 8937996
$ws.Range("I61").Value = 11911679
$ws.Range("J61").Value = 16946.785
$ws.Range("K61").Value = 11911679
$ws.Range("L61").Value = 16946.785
$ws.Range("M61").Value = -11911467
$ws.Range("N61").Value = -17370.785

# Row 74
$ws.Range("H74").Value = 4284.5347
$ws.Range("I74").Value = 3470.4595
$ws.Range("K74").Value = 3470.4595
$ws.Range("M74").Value = -2596.4595

# Row 77
$ws.Range("H77").Value = 4284.5347
$ws.Range("I77").Value = 3470.4595
$ws.Range("K77").Value = 17352.2975
$ws.Range("M77").Value = -12984.2975

# Row 97
$ws.Range("H97").Value = 1103.4667
$ws.Range("I97").Value = 1050.24
$ws.Range("J97").Value = 1369.6
$ws.Range("K97").Value = 1050.24
$ws.Range("L97").Value = 1369.6
$ws.Range("M97").Value = -554.24
$ws.Range("N97").Value = -2361.6

# Row 122
$ws.Range("H122").Value = 1430.08
$ws.Range("I122").Value = 1378.3684
$ws.Range("J122").Value = 1593.8334
$ws.Range("K122").Value = 4135.1052
$ws.Range("L122").Value = 4781.5002
$ws.Range("M122").Value = -1685.1052
$ws.Range("N122").Value = -9681.5002

# Row 124
$ws.Range("H124").Value = 35985.8
$ws.Range("J124").Value = 35985.8
$ws.Range("L124").Value = 35985.8
$ws.Range("N124").Value = -45805.8

# Row 136
$ws.Range("H136").Value = 8937996
$ws.Range("I136").Value = 11911679
$ws.Range("J136").Value = 16946.785
$ws.Range("K136").Value = 35735037
$ws.Range("L136").Value = 50840.355
$ws.Range("M136").Value = -35732487
$ws.Range("N136").Value = -55940.355

# Row 139
$ws.Range("H139").Value = 119665
$ws.Range("J139").Value = 119665
$ws.Range("L139").Value = 119665
$ws.Range("N139").Value = -129945

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2501.6155
$ws.Range("I20").Value = 2300.8076
$ws.Range("J20").Value = 2903.2307
$ws.Range("K20").Value = 2300.8076
$ws.Range("L20").Value = 2903.2307
$ws.Range("M20").Value = -2053.8076
$ws.Range("N20").Value = -3397.2307

# Row 110
$ws.Range("H110").Value = 27666.334
$ws.Range("J110").Value = 27666.334
$ws.Range("L110").Value = 27666.334
$ws.Range("N110").Value = -35846.334

# Row 111
$ws.Range("H111").Value = 31249.5
$ws.Range("J111").Value = 31249.5
$ws.Range("L111").Value = 31249.5
$ws.Range("N111").Value = -39429.5

# Row 134
$ws.Range("H134").Value = 3420.6
$ws.Range("I134").Value = 3566.1155
$ws.Range("J134").Value = 898.3333
$ws.Range("K134").Value = 10698.3465
$ws.Range("L134").Value = 2694.9999
$ws.Range("M134").Value = -8163.3465
$ws.Range("N134").Value = -7764.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4582.775
$ws.Range("I31").Value = 2142.2104
$ws.Range("K31").Value = 2142.2104
$ws.Range("M31").Value = -1847.2104

# Row 34
$ws.Range("H34").Value = 4582.775
$ws.Range("I34").Value = 2142.2104
$ws.Range("K34").Value = 2142.2104
$ws.Range("M34").Value = -1940.2104

# Row 58
$ws.Range("H58").Value = 6973.9287
$ws.Range("I58").Value = 5490.067
$ws.Range("J58").Value = 8686.076999999999
$ws.Range("K58").Value = 5490.067
$ws.Range("L58").Value = 8686.076999999999
$ws.Range("M58").Value = -5287.067
$ws.Range("N58").Value = -9092.076999999999

# Row 62
$ws.Range("H62").Value = 6199.3887
$ws.Range("I62").Value = 5860.4443
$ws.Range("J62").Value = 6538.3335
$ws.Range("K62").Value = 5860.4443
$ws.Range("L62").Value = 6538.3335
$ws.Range("M62").Value = -5236.4443
$ws.Range("N62").Value = -7786.3335

# Row 65
$ws.Range("H65").Value = 6199.3887
$ws.Range("I65").Value = 5860.4443
$ws.Range("J65").Value = 6538.3335
$ws.Range("K65").Value = 29302.2215
$ws.Range("L65").Value = 32691.6675
$ws.Range("M65").Value = -26182.2215
$ws.Range("N65").Value = -38931.6675

# Row 86
$ws.Range("H86").Value = 51914.418
$ws.Range("I86").Value = 10662.167
$ws.Range("K86").Value = 10662.167
$ws.Range("M86").Value = -9539.166999999999

# Row 89
$ws.Range("H89").Value = 51914.418
$ws.Range("I89").Value = 10662.167
$ws.Range("K89").Value = 53310.835
$ws.Range("M89").Value = -47694.835

# Row 132
$ws.Range("H132").Value = 3388.3333
$ws.Range("I132").Value = 3459.673
$ws.Range("J132").Value = 2924.625
$ws.Range("K132").Value = 10379.019
$ws.Range("L132").Value = 8773.875
$ws.Range("M132").Value = -7849.019
$ws.Range("N132").Value = -13833.875

# Row 134
$ws.Range("H134").Value = 4791.1562
$ws.Range("I134").Value = 3724.24
$ws.Range("K134").Value = 11172.72
$ws.Range("M134").Value = -8637.719999999999

# Row 136
$ws.Range("H136").Value = 6973.9287
$ws.Range("I136").Value = 5490.067
$ws.Range("J136").Value = 8686.076999999999
$ws.Range("K136").Value = 16470.201
$ws.Range("L136").Value = 26058.231
$ws.Range("M136").Value = -13920.201
$ws.Range("N136").Value = -31158.231

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 300
$ws.Range("L13").ClearContents()
$ws.Range("M13").Value = -132
$ws.Range("N13").Value = 0

# Row 113
$ws.Range("H113").Value = 1038
$ws.Range("I113").Value = 317
$ws.Range("K113").Value = 951
$ws.Range("M113").Value = 1219

# Row 123
$ws.Range("H123").Value = 14777.375
$ws.Range("I123").Value = 3407.3333
$ws.Range("K123").Value = 10221.9999
$ws.Range("M123").Value = -7771.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 48750
$ws.Range("J5").Value = 48750
$ws.Range("L5").Value = 48750
$ws.Range("N5").Value = -48974

# Row 113
$ws.Range("H113").Value = 2799.0908
$ws.Range("I113").Value = 1865.5555
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 1865.5555
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 304.4445000000001
$ws.Range("N113").Value = -11340

# Row 122
$ws.Range("H122").Value = 6716.7617
$ws.Range("J122").Value = 13198.5
$ws.Range("L122").Value = 39595.5
$ws.Range("N122").Value = -44495.5

# Row 136
$ws.Range("H136").Value = 55897.785
$ws.Range("J136").Value = 55897.785
$ws.Range("L136").Value = 167693.355
$ws.Range("N136").Value = -172793.355

# Row 140
$ws.Range("H140").Value = 140000
$ws.Range("J140").Value = 140000
$ws.Range("L140").Value = 140000
$ws.Range("N140").Value = -150360

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1666.3334
$ws.Range("I22").Value = 1764.125
$ws.Range("K22").Value = 1764.125
$ws.Range("M22").Value = -1469.125

# Row 27
$ws.Range("H27").Value = 1666.3334
$ws.Range("I27").Value = 1764.125
$ws.Range("K27").Value = 1764.125
$ws.Range("M27").Value = -1657.125

# Row 40
$ws.Range("H40").Value = 4130.077
$ws.Range("I40").Value = 3586.375
$ws.Range("K40").Value = 3586.375
$ws.Range("M40").Value = -3450.375

# Row 42
$ws.Range("H42").Value = 25235.273
$ws.Range("J42").Value = 25411.334
$ws.Range("L42").Value = 25411.334
$ws.Range("N42").Value = -26537.334

# Row 49
$ws.Range("H49").Value = 25235.273
$ws.Range("J49").Value = 25411.334
$ws.Range("L49").Value = 25411.334
$ws.Range("N49").Value = -25705.334

# Row 55
$ws.Range("H55").Value = 391.06668
$ws.Range("I55").Value = 504.36365
$ws.Range("J55").Value = 79.5
$ws.Range("K55").Value = 504.36365
$ws.Range("L55").Value = 79.5
$ws.Range("M55").Value = -331.36365
$ws.Range("N55").Value = -425.5

# Row 56
$ws.Range("H56").Value = 21990
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 93
$ws.Range("H93").Value = 5842.2246
$ws.Range("I93").Value = 1503.3143
$ws.Range("K93").Value = 1503.3143
$ws.Range("M93").Value = -255.3143

# Row 132
$ws.Range("H132").Value = 12578.754
$ws.Range("I132").Value = 11546.967
$ws.Range("K132").Value = 34640.901
$ws.Range("M132").Value = -32110.901

# Row 136
$ws.Range("H136").Value = 4234.1665
$ws.Range("I136").Value = 4810.5557
$ws.Range("K136").Value = 14431.6671
$ws.Range("M136").Value = -11881.6671

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3008.8057
$ws.Range("I122").Value = 2214.9285
$ws.Range("J122").Value = 5787.375
$ws.Range("K122").Value = 6644.7855
$ws.Range("L122").Value = 17362.125
$ws.Range("M122").Value = -4194.7855
$ws.Range("N122").Value = -22262.125

# Row 132
$ws.Range("H132").Value = 3123.6956
$ws.Range("I132").Value = 3288.1628
$ws.Range("K132").Value = 9864.4884
$ws.Range("M132").Value = -7334.4884

# Row 133
$ws.Range("H133").Value = 145178.25
$ws.Range("J133").Value = 145178.25
$ws.Range("L133").Value = 145178.25
$ws.Range("N133").Value = -155298.25

# Row 136
$ws.Range("H136").Value = 4330.816
$ws.Range("I136").Value = 2201.7856
$ws.Range("J136").Value = 10292.1
$ws.Range("K136").Value = 6605.3568
$ws.Range("L136").Value = 30876.3
$ws.Range("M136").Value = -4055.3568
$ws.Range("N136").Value = -35976.3
